$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username/password test data in rows 2 and 3
$ws.Range("B2").Value = "shopfloor1"
$ws.Range("C2").Value = "shopfloor1*1"
$ws.Range("B3").Value = "shopfloor1"
$ws.Range("C3").Value = "shopfloor1*1"

# Update the selected range as reflected in the saved sheet view
$ws.Range("B3:C3").Select()
